$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2-57 is being bumped from
# serial date 45179 (2023-09-10) to 45180 (2023-09-11).
for ($row = 2; $row -le 57; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value = 45180
    }
}
